$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.902.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +5.40%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.355.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +5.42%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'570.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +6.88%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'152.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +5.87%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.09%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.357.56"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +5.46%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.528"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.69%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'7.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.54%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +5.12%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.440"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.53%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.925.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +5.17%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -0.04%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'26.92"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +4.10%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0000180"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +4.75%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'WrappedEther"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'3.440.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +7.90%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "'WrappedBTC"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'62.876.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +5.37%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +1.69%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'13.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +5.83%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'8.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +2.81%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'385.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +5.23%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.11%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +2.53%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'70.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.51%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'9.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +6.57%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +6.72%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0₃0963"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +10.11%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.14%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +6.21%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'22.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +3.10%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'5.55"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +5.02%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +9.57%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'6.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +3.10%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'6.68"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.23%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +9.31%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'157.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.26%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.87"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +12.24%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'26.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +3.74%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +12.67%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.0739"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +6.06%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'2.842.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.03%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +3.82%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'4.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.35%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.743"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +4.93%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'ONDO"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'1.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +5.14%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'RenzoRestakedETH"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'3.394.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +5.36%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'21.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +7.34%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'Bittensor"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'296.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +13.24%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'Stellar"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'0.103"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.84%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +2.76%  "
$ws.Range("E51").Style = "Normal"
